$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3-data")

# Rename the slug/title values for the widget-menu related rows (t2x -> t3x, suffix "3")
# Order matters for shared-string table append order, matching the target workbook.
$ws.Range("A3").Value = "t31"
$ws.Range("A5").Value = "t32"
$ws.Range("C3").Value = "basic_geo_db3"
$ws.Range("C4").Value = "nature_geo3"
$ws.Range("C5").Value = "nature_resources3"
$ws.Range("C6").Value = "nature_res3"
$ws.Range("C7").Value = "eco_res3"

# Update the active selection on the sheet
$ws.Range("C8").Select()
